$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D can look like numbers (or dates) to Excel's
# auto-detection (e.g. "1.002", "6.611"). Force text storage by setting
# NumberFormat to "@" before the write, then restore the default "Normal"
# style afterwards so no stray number-format style is left on the cell.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.428.13"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "2.018.36"
$ws.Range("E3").Value = "  +5.93%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "324.95"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Value = "0.5131"
$ws.Range("E7").Value = "  +2.11%  "
$ws.Range("D8").Value = "0.4223"
$ws.Range("E8").Value = "  +4.46%  "
$ws.Range("D9").Value = "0.08717"
$ws.Range("E9").Value = "  +5.61%  "
$ws.Range("E10").Value = "  +3.85%  "
$ws.Range("D11").Value = "1.137"
$ws.Range("E11").Value = "  +3.84%  "
$ws.Range("D12").Value = "24.87"
$ws.Range("E12").Value = "  +3.72%  "
$ws.Range("D13").Value = "2.018.00"
$ws.Range("E13").Value = "  +5.87%  "
$ws.Range("D14").Value = "6.611"
$ws.Range("E14").Value = "  +4.08%  "
$ws.Range("D15").Value = "7.490"
$ws.Range("E15").Value = "  +4.22%  "
$ws.Range("E16").Value = "  -0.18%  "
$ws.Range("D17").Value = "94.48"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").Value = "0.00001115"
$ws.Range("E18").Value = "  +2.14%  "
$ws.Range("D19").Value = "0.06532"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "18.98"
$ws.Range("E20").Value = "  +5.42%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "6.212"
$ws.Range("E22").Value = "  +4.77%  "
$ws.Range("D23").Value = "30.484.33"
$ws.Range("E23").Value = "  +1.72%  "
$ws.Range("E24").Value = "  +5.74%  "
$ws.Range("D25").Value = "2.228"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "2.254.95"
$ws.Range("E26").Value = "  +6.10%  "
$ws.Range("D27").Value = "22.46"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("D28").Value = "162.76"
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "2.428"
$ws.Range("E29").Value = "  +7.25%  "
$ws.Range("D30").Value = "131.54"
$ws.Range("E30").Value = "  +2.35%  "
$ws.Range("D31").Value = "1.143"
$ws.Range("E31").Value = "  +2.02%  "
$ws.Range("D32").Value = "0.1053"
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("D33").Value = "6.073"
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("D35").Value = "1.370"
$ws.Range("E35").Value = "  +15.20%  "
$ws.Range("D36").Value = "0.02531"
$ws.Range("E36").Value = "  +4.18%  "
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "0.06674"
$ws.Range("E38").Value = "  +5.50%  "
$ws.Range("D39").Value = "12.30"
$ws.Range("E39").Value = "  +8.93%  "
$ws.Range("D40").Value = "0.2201"
$ws.Range("E40").Value = "  +3.03%  "
$ws.Range("D41").Value = "9.076"
$ws.Range("E41").Value = "  +5.41%  "
$ws.Range("D42").Value = "0.6676"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").Value = "1.232"
$ws.Range("E43").Value = "  +2.68%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "13.68"
$ws.Range("E45").Value = "  +3.00%  "
$ws.Range("D46").Value = "0.6191"
$ws.Range("E46").Value = "  +3.15%  "
$ws.Range("D47").Value = "2.192"
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").Value = "3.659"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("E49").Value = "  +5.36%  "
$ws.Range("D50").Value = "124.75"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").Value = "81.13"
$ws.Range("E51").Value = "  +3.62%  "

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}
